# Apply the "feature/scorpionssuite-refactoring" VIN upload-test fixture update:
# append 4 new sample data rows (6-9) that mirror rows 2-5 but use the new
# SYMBOL_2017 version label and unique BI/PD/UM/MP symbol codes per row,
# then move the active selection to J17 to match the saved worksheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column numbers (within row 6-9 new data) that reuse the worksheet's existing
# "left aligned" cell style (style index 2 in the original file) - this mirrors
# every column in rows 2-5 except A, B, D, E, F, J, L, which are left unstyled.
$styledColumns = @(3, 7, 8, 9, 11, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36)

# New rows 6-9, one array of 36 values (columns A..AJ) per row.
$newRows = @(
    @("GGGKN3DD&E", "SYMBOL_2017", 2018, "TOYOTA", "TOYOTA", "Gt", "MDX ADVANCE", 53080, "WAG", "UT_SS", "SUV", "UT_SS", "WAG", "4.5L V10", 8, "G", 214, "2WD", 2, "000R", "DUAL AIR BAGS FRONT", 2, "4 WHEEL STANDARD", "STD", "B-IMMOBILIZER/KEYLSS ENTRY/ALARM", 42, 42, "Y", "BI001", "PD001", "UM001", "MP001", 20010101, "Y", "Y", "N"),
    @("GGGKN3DD&E", "SYMBOL_2017", 2018, "UT_SS", "UT_SS", "invalidVIN", "MDX ADVANCE", 53080, "WAG", "UT_SS", "SUV", "UT_SS", "WAG", "4.5L V10", 8, "G", 214, "2WD", 2, "000R", "DUAL AIR BAGS FRONT", 2, "4 WHEEL STANDARD", "STD", "B-IMMOBILIZER/KEYLSS ENTRY/ALARM", 42, 42, "Y", "BI002", "PD002", "UM002", "MP002", 20000101, "N", "Y", "N"),
    @("GGGKN3DD&E", "SYMBOL_2017", 2018, "UT_SS", "UT_SS", "SecondValid", "MDX ADVANCE", 53080, "WAG", "UT_SS", "SUV", "UT_SS", "WAG", "4.5L V10", 8, "G", 214, "2WD", 2, "000R", "DUAL AIR BAGS FRONT", 2, "4 WHEEL STANDARD", "STD", "B-IMMOBILIZER/KEYLSS ENTRY/ALARM", 42, 42, "Y", "BI003", "PD003", "UM003", "MP003", 20150101, "Y", "Y", "N"),
    @("GGGKN3DD&E", "SYMBOL_2017", 2018, "UT_SS", "UT_SS", "ThirdValid", "MDX ADVANCE", 53080, "WAG", "UT_SS", "SUV", "UT_SS", "WAG", "4.5L V10", 8, "G", 214, "2WD", 2, "000R", "DUAL AIR BAGS FRONT", 2, "4 WHEEL STANDARD", "STD", "B-IMMOBILIZER/KEYLSS ENTRY/ALARM", 42, 42, "Y", "BI004", "PD004", "UM004", "MP004", 20190101, "Y", "Y", "N")
)

$startRow = 6
foreach ($rowValues in $newRows) {
    $colIndex = 1
    foreach ($value in $rowValues) {
        if ($styledColumns -contains $colIndex) {
            $ws.Cells.Item($startRow, $colIndex).HorizontalAlignment = -4131
        }
        $ws.Cells.Item($startRow, $colIndex).Value = $value
        $colIndex = $colIndex + 1
    }
    $startRow = $startRow + 1
}

# Match the saved sheet view: selection moved to J17 (dimension auto-extends
# to A1:AJ9 once the new rows above are written).
$ws.Range("J17").Select()
